$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values are text that may look numeric (e.g. "1.014").
# Force text format so Excel does not auto-convert them to numbers,
# then restore the default "Normal" style so no stray formatting is left behind.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.762.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.07%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.852.89'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.60%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.014'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -2.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '320.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.24%  '

$ws.Range("E6").Value = '  -2.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4323'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.08%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3767'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.75%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07395'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.90%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8855'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.14%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.72'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.860.45'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.767'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.483'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.36%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.07152'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.70%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.45%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.015'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.25%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009047'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.99%  '

$ws.Range("E19").Value = '  -2.29%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.54'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.726.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.29%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.279'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.71%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.22'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.105.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.63%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.025'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.64'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.140'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.40%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.438'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.16%  '

$ws.Range("E30").Value = '  +2.77%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08967'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.57%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.242'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.91%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7847'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.96%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.584'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.924'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.88%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.147'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.74%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.013'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.37%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05342'
$ws.Range("D38").Style = "Normal"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01971'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.22%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.151'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.869'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.73%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5196'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.17%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1689'
$ws.Range("D43").Style = "Normal"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '9.129'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.66%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '110.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.74'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.09%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.725'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.4751'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.79%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06518'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.08%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.013'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.47%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.906'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.25%  '

